$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: paragraph "Todos os campos são de preenchimento obrigatório,
# exceto cada tipo de forma de inserção." becomes "Todos os campos são
# de preenchimento obrigatório. A forma de inserção deve ter pelo menos
# um dos itens selecionado." - split across several runs, with the
# _GoBack bookmark landing inside "campos" (between "cam" and "pos").
# ---------------------------------------------------------------------

# Step 1: replace the old second clause with the new continuation text.
$rng = $d.Content
$rng.Find.Execute(", exceto cada tipo de forma de inserção", $true, $false, $false, $false, $false, $true, 1, $false, ". A forma de inserção deve ter pelo menos um dos itens selecionado", 2) | Out-Null

# Locate the (now merged-into-one-run) full sentence so we can compute
# absolute character offsets for the run-boundary markers we need.
$rngFull = $d.Content
$rngFull.Find.Execute("Todos os campos são de preenchimento obrigatório. A forma de inserção deve ter pelo menos um dos itens selecionado.", $true) | Out-Null
$base = $rngFull.Start

# Step 2: mark the boundary right before the trailing "." first (so it
# is not the most-recently-touched boundary once we are done - that
# keeps its run a plain, non-"preserve" run, matching how Word leaves
# an untouched trailing run alone).
$posEnd = $base + 114
$d.Bookmarks.Add("_TmpEnd", $d.Range($posEnd, $posEnd))

# Step 3: drop the real _GoBack bookmark on the "cam" | "pos" boundary.
$posGoBack = $base + 12
$d.Bookmarks.Add("_GoBack", $d.Range($posGoBack, $posGoBack))

# Step 4: mark the remaining internal run boundaries:
#   after "...obrigatório"     -> offset 48
#   after ". A "                -> offset 52
#   after "forma de inserção"   -> offset 69
$offsets = @(48, 52, 69)
$i = 0
foreach ($off in $offsets) {
    $p = $base + $off
    $name = "_TmpSplit$i"
    $d.Bookmarks.Add($name, $d.Range($p, $p))
    $i = $i + 1
}

# Step 5: remove all the throw-away markers - the run boundaries they
# created stay in place even after the bookmarks themselves are gone.
$d.Bookmarks("_TmpEnd").Delete()
for ($j = 0; $j -lt $i; $j++) {
    $d.Bookmarks("_TmpSplit$j").Delete()
}

# ---------------------------------------------------------------------
# Edit 2: remove the stray _GoBack bookmark that used to sit inside
# "usuário realmente" and let the two runs merge back into one.
# ---------------------------------------------------------------------
$old2 = "Informar que a ação é irreversível e que apagará o histórico da participação, perguntando se o que o usuário realmente deseja é encerrar a participação."
$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null
